$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 6 (LENS) - update Giornata 10 through Giornata 18 cumulative values (columns K:S)
$ws.Range("K6").Value = 7.455798300000001
$ws.Range("L6").Value = 8.485788299999999
$ws.Range("M6").Value = 9.1578713
$ws.Range("N6").Value = 11.7852413
$ws.Range("O6").Value = 12.2890633
$ws.Range("P6").Value = 13.9062533
$ws.Range("Q6").Value = 14.8347443
$ws.Range("R6").Value = 15.9951243
$ws.Range("S6").Value = 17.1716043

# Row 10 (MARSIGLIA) - update Giornata 10 through Giornata 18 cumulative values (columns K:S)
$ws.Range("K10").Value = 14.027083
$ws.Range("L10").Value = 14.351284
$ws.Range("M10").Value = 14.351284
$ws.Range("N10").Value = 14.673567
$ws.Range("O10").Value = 15.0582
$ws.Range("P10").Value = 16.76591
$ws.Range("Q10").Value = 17.9344
$ws.Range("R10").Value = 18.297546
$ws.Range("S10").Value = 19.077934
